$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# New column F header: "schema:startDate^^xsd:date", same header style as the other
# header cells on row 7 (copy format from E7 so the existing style index is reused).
$ws.Range("F7").Value = "schema:startDate^^xsd:date"
$ws.Range("E7").Copy()
$ws.Range("F7").PasteSpecial(-4122)

# New date values in column F for the two data rows, formatted as a short date
# (numFmtId 14) like Excel would apply automatically for a date value. Format F8
# first, then copy its formatting onto F9 so both cells share a single style record
# (instead of each NumberFormat assignment minting its own duplicate style).
$ws.Range("F8").Value = 42695
$ws.Range("F8").NumberFormat = "mm-dd-yy"

$ws.Range("F9").Value = 42519
$ws.Range("F8").Copy()
$ws.Range("F9").PasteSpecial(-4122)

# Give the new column a sensible custom width (matches the author's width as closely
# as the host's column-width model allows).
$ws.Columns.Item(6).ColumnWidth = 30.7369791666

# Move the active selection to follow the newly-added column.
[void]$ws.Activate()
[void]$ws.Range("F10").Select()

# Unrelated formatting fix on the second sheet: the header row no longer needs an
# explicit (taller) row height, so let it fall back to the sheet's default.
$ws2 = $wb.Worksheets.Item("Locations")
[void]$ws2.Rows.Item(6).AutoFit()
